# "fixed land cover on the map"
#
# 1) Bump the cached auto-date footer field (datetimeFigureOut) on the
#    slide master and on every slide layout from 10/8/2024 -> 10/10/2024.
# 2) Fix the mislabeled land-cover percentages on the map legend
#    (Group 14 inside the map group on slide 1).

$p = $ppt.ActivePresentation

# --- 1. Update the cached date placeholder text wherever it appears ----
$oldDate = "10/8/2024"
$newDate = "10/10/2024"

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

$layouts = $master.CustomLayouts
for ($l = 1; $l -le $layouts.Count; $l++) {
    $layout = $layouts.Item($l)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# --- 2. Fix the land-cover legend percentages on the map (slide 1) -----
$s = $p.Slides.Item(1)
$mapGroup = $s.Shapes.Item(1)
$items = $mapGroup.GroupItems

# Corrected value per legend textbox, addressed by shape Name (not by its
# old text) since several of the old values overlap with the values other
# boxes are being changed to (e.g. 27.3% / 27.6% / 14.8% swap around).
$newValues = @{
    "TextBox 7"  = "1.25%"
    "TextBox 8"  = "27.6%"
    "TextBox 10" = "14.8%"
    "TextBox 11" = "5.17%"
    "TextBox 12" = "27.3%"
}

for ($i = 1; $i -le $items.Count; $i++) {
    $shp = $items.Item($i)
    if ($newValues.ContainsKey($shp.Name)) {
        $shp.TextFrame.TextRange.Text = $newValues[$shp.Name]
    }
}
